$wb = $excel.ActiveWorkbook

# Reference the existing "总计" (Total) sheet, used as insertion anchor
$totalSheetAnchor = $wb.Worksheets.Item("总计")

# Insert the new "2022-Q1" worksheet right before "总计" (matches final sheet order)
$newWs = $wb.Worksheets.Add($totalSheetAnchor, $null)
$newWs.Name = "2022-Q1"

# Re-fetch the sheet objects by name now that a new sheet has been inserted,
# since previously captured references can become stale/reindexed.
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet = $wb.Worksheets.Item("2021-Q4")
$newWs = $wb.Worksheets.Item("2022-Q1")

# ---- Header row (B1:H1) ----
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# ---- Data rows ----
# Columns B (fund codes w/ leading zeros) and D,E,F,G (stored as text, matching
# the formatting convention used on the other quarterly sheets) must be forced
# to text format before the values are assigned, otherwise Excel will coerce
# them to numbers and leading zeros / trailing zeros would be lost.
$newWs.Range("B2:B7").NumberFormat = "@"
$newWs.Range("D2:G7").NumberFormat = "@"

$newWs.Range("A2").Value = 0
$newWs.Range("B2").Value = "011046"
$newWs.Range("C2").Value = "富国优质企业混合A"
$newWs.Range("D2").Value = "8.18"
$newWs.Range("E2").Value = "71.23"
$newWs.Range("F2").Value = "2.83"
$newWs.Range("G2").Value = "0.2315"
$newWs.Range("H2").Value = 7

$newWs.Range("A3").Value = 1
$newWs.Range("B3").Value = "009782"
$newWs.Range("C3").Value = "富国兴泉回报12个月持有期混合A"
$newWs.Range("D3").Value = "6.18"
$newWs.Range("E3").Value = "70.06"
$newWs.Range("F3").Value = "2.22"
$newWs.Range("G3").Value = "0.1372"
$newWs.Range("H3").Value = 10

$newWs.Range("A4").Value = 2
$newWs.Range("B4").Value = "012096"
$newWs.Range("C4").Value = "鑫元鑫动力混合型证券投资基金A"
$newWs.Range("D4").Value = "2.82"
$newWs.Range("E4").Value = "88.33"
$newWs.Range("F4").Value = "3.77"
$newWs.Range("G4").Value = "0.1063"
$newWs.Range("H4").Value = 9

$newWs.Range("A5").Value = 3
$newWs.Range("B5").Value = "009783"
$newWs.Range("C5").Value = "富国兴泉回报12个月持有期混合C"
$newWs.Range("D5").Value = "2.45"
$newWs.Range("E5").Value = "70.06"
$newWs.Range("F5").Value = "2.22"
$newWs.Range("G5").Value = "0.0544"
$newWs.Range("H5").Value = 10

$newWs.Range("A6").Value = 4
$newWs.Range("B6").Value = "005732"
$newWs.Range("C6").Value = "富国臻选成长灵活配置混合"
$newWs.Range("D6").Value = "2.45"
$newWs.Range("E6").Value = "64.81"
$newWs.Range("F6").Value = "2.15"
$newWs.Range("G6").Value = "0.0527"
$newWs.Range("H6").Value = 10

$newWs.Range("A7").Value = 5
$newWs.Range("B7").Value = "011047"
$newWs.Range("C7").Value = "富国优质企业混合C"
$newWs.Range("D7").Value = "0.48"
$newWs.Range("E7").Value = "71.23"
$newWs.Range("F7").Value = "2.83"
$newWs.Range("G7").Value = "0.0136"
$newWs.Range("H7").Value = 7

# ---- Apply formatting to match the style used on the other quarterly sheets ----
# Header row B1:H1 and the index column A2:A7 use bold-centered style (s="2")
$templateSheet.Range("B1:H1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)

$templateSheet.Range("A2:A7").Copy()
$newWs.Range("A2:A7").PasteSpecial(-4122)

# ---- Update the "总计" (Total) sheet: insert a new top data row for 2022-Q1 ----
$totalSheet.Rows.Item(2).Insert()

# Excel's row-Insert carries formatting down from the header row onto the new
# row; clear that back to the default (unstyled) look used by the other data
# rows in this sheet (columns B:D have no explicit style there).
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 6
$totalSheet.Range("D2").Value = 0.6

# Re-apply the index-column style to the newly inserted A2 cell
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# Renumber the index column (A) for the rows pushed down
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
